$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(10114.89, 10201.6, 18.84, 19, $true, 0.85, 42613.766944444447, $false),
    @(10093.65, 10114.89, 18.93, 18.97, $true, 0.21, 42614.675393518519, $false),
    @(10001.799999999999, 10093.65, 18.72, 18.89, $true, 0.91, 42615.752129629633, $false)
)

$r = 7
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]

    # Reuse the existing date/time number format (style index already used by
    # column G in rows 3-6) instead of letting NumberFormat mint a brand new
    # style entry - copy formatting from row 3's date cell.
    $ws.Cells.Item(3, 7).Copy()
    $ws.Cells.Item($r, 7).PasteSpecial(-4122)

    $r = $r + 1
}

$excel.CutCopyMode = 0
